$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 previously held the "is_locked_lbl" dictionary-select template string.
# The is_locked / is_enabled columns are removed, so the remaining columns
# shift left: D1 now carries what used to be the "order_by" column text.
$ws.Range("D1").Value = '<%=comment.order_by%>'

# E1 previously held the "is_enabled_lbl" dictionary-select template string.
# It now carries what used to be the "rem" column text.
$ws.Range("E1").Value = '<%=comment.rem%>'

# F1 previously held the "order_by" text; the new tenant_id column is
# appended, so F1 now carries the new tenant_id_lbl select-list template.
$ws.Range("F1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'

# G1 (old "rem" column) is no longer part of the row; remove it entirely.
$ws.Range("G1").ClearContents()
